# edit.ps1 - apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17 and 18 swapped positions (Polkadot moved above WrappedEther)
$ws.Range("B17").Value2 = "Polkadot"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "5.90"
$ws.Range("E17").Value2 = "  -5.35%  "

$ws.Range("B18").Value2 = "WrappedEther"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value2 = "2.932.82"
$ws.Range("E18").Value2 = "  -4.05%  "

# Updated Price (D) and Volume(1h) (E) values for remaining rows
$ws.Range("D2").Value2 = "56.002.37"
$ws.Range("E2").Value2 = "  -3.60%  "
$ws.Range("D3").Value2 = "2.949.15"
$ws.Range("E3").Value2 = "  -3.52%  "
$ws.Range("E4").Value2 = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "491.87"
$ws.Range("E5").Value2 = "  -6.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "132.89"
$ws.Range("E6").Value2 = "  -6.56%  "
$ws.Range("E7").Value2 = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.420"
$ws.Range("E8").Value2 = "  -6.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "7.09"
$ws.Range("E9").Value2 = "  -6.02%  "
$ws.Range("E10").Value2 = "  -7.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.347"
$ws.Range("E11").Value2 = "  -5.93%  "
$ws.Range("D12").Value2 = "3.442.35"
$ws.Range("E12").Value2 = "  -3.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.124"
$ws.Range("E13").Value2 = "  -4.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "25.84"
$ws.Range("E14").Value2 = "  -5.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.0000156"
$ws.Range("E15").Value2 = "  -9.33%  "
$ws.Range("D16").Value2 = "56.047.38"
$ws.Range("E16").Value2 = "  -3.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "12.41"
$ws.Range("E19").Value2 = "  -5.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "7.69"
$ws.Range("E20").Value2 = "  -5.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "310.97"
$ws.Range("E21").Value2 = "  -8.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "1.00"
$ws.Range("E22").Value2 = "  +0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "5.76"
$ws.Range("E23").Value2 = "  +0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.480"
$ws.Range("E24").Value2 = "  -4.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "62.21"
$ws.Range("E25").Value2 = "  -4.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.994"
$ws.Range("E26").Value2 = "  -0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.160"
$ws.Range("E27").Value2 = "  -5.23%  "
$ws.Range("D28").Value2 = "0.0₃0834"
$ws.Range("E28").Value2 = "  -13.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "6.37"
$ws.Range("E29").Value2 = "  -8.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "6.96"
$ws.Range("E30").Value2 = "  -8.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "1.73"
$ws.Range("E31").Value2 = "  -6.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "19.84"
$ws.Range("E32").Value2 = "  -6.25%  "
$ws.Range("E33").Value2 = "  -10.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "149.35"
$ws.Range("E34").Value2 = "  -6.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "4.41"
$ws.Range("E35").Value2 = "  -7.90%  "
$ws.Range("E36").Value2 = "  -6.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "1.19"
$ws.Range("E37").Value2 = "  -9.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "23.77"
$ws.Range("E38").Value2 = "  -7.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.0648"
$ws.Range("E39").Value2 = "  -6.97%  "
$ws.Range("D40").Value2 = "2.977.82"
$ws.Range("E40").Value2 = "  -3.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "36.46"
$ws.Range("E41").Value2 = "  -3.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.999"
$ws.Range("E42").Value2 = "  -0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "3.64"
$ws.Range("E43").Value2 = "  -7.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.628"
$ws.Range("E44").Value2 = "  -5.84%  "
$ws.Range("D45").Value2 = "2.119.76"
$ws.Range("E45").Value2 = "  -9.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "1.33"
$ws.Range("E46").Value2 = "  -9.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "5.81"
$ws.Range("E47").Value2 = "  -4.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.902"
$ws.Range("E48").Value2 = "  -12.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0228"
$ws.Range("E49").Value2 = "  -6.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "18.70"
$ws.Range("E50").Value2 = "  -6.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.0838"
$ws.Range("E51").Value2 = "  -7.33%  "
